$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the sheet
$ws.Name = "LoginInformation"

# Update row 2 values (C2, D2) - the batch credentials
$ws.Range("C2").Value = "c"
$ws.Range("D2").Value = "B"

# Add new row 3 with TradBatch login info
$ws.Range("A3").Value = "TradBatch"
$ws.Range("B3").Value = "http://gb2trpec-001.ffastserve.com/j2ee/"
$ws.Range("C3").Value = "a"
$ws.Range("D3").Value = "b"

# Copy formatting from row 2 down to row 3 so it reuses the existing bordered/text style
$src = $ws.Range("A2:D2")
$dst = $ws.Range("A3:D3")
$src.Copy() | Out-Null
$dst.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Widen column B to fit the longer URL text (closest achievable value to the
# target bestFit width of 38.7109375 given this runtime's column-width rounding)
$ws.Columns.Item(2).ColumnWidth = 38.0

# Move the active selection to A4, as it was left after entering the new row
$ws.Activate() | Out-Null
$ws.Range("A4").Select() | Out-Null
